$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) updates: force text entry so number-like strings
# (e.g. "575.55") are stored as text, matching the source inlineStr cells,
# instead of being auto-parsed into floating point numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '61.564.56'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.390.11'
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.25'
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.67'
$ws.Range("D9").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.968.23'
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.44'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.400.63'
$ws.Range("D15").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.591.16'
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.58'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.98'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '390.72'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.04'
$ws.Range("D22").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000113'
$ws.Range("D25").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.40'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.00'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.27'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.89'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '168.71'
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.424.63'
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0763'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.87'
$ws.Range("D40").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.475.99'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.85'
$ws.Range("D46").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.04'
$ws.Range("D50").Style = "Normal"

# Coin name / link / volume column updates (plain text already).
$ws.Range("E2").Value = '  +0.71%  '
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("E6").Value = '  -0.87%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  -0.78%  '
$ws.Range("E9").Value = '  +0.94%  '
$ws.Range("E11").Value = '  -3.07%  '
$ws.Range("E12").Value = '  -0.38%  '
$ws.Range("E13").Value = '  +0.21%  '
$ws.Range("E14").Value = '  +1.29%  '
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("E16").Value = '  -0.91%  '
$ws.Range("E17").Value = '  +0.69%  '
$ws.Range("E18").Value = '  -0.28%  '
$ws.Range("E19").Value = '  -2.11%  '
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("E21").Value = '  +1.80%  '
$ws.Range("E22").Value = '  +0.70%  '
$ws.Range("E23").Value = '  -1.51%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("E25").Value = '  -4.68%  '
$ws.Range("E26").Value = '  +7.83%  '
$ws.Range("E27").Value = '  -0.06%  '
$ws.Range("E28").Value = '  -1.61%  '
$ws.Range("E30").Value = '  -1.14%  '
$ws.Range("B31").Value = 'Fetch.AI'
$ws.Range("C31").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("E31").Value = '  -1.32%  '
$ws.Range("B32").Value = 'USDe'
$ws.Range("C32").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("E32").Value = '  +0.04%  '
$ws.Range("E33").Value = '  -1.06%  '
$ws.Range("E34").Value = '  -1.82%  '
$ws.Range("E35").Value = '  +0.57%  '
$ws.Range("E36").Value = '  -0.02%  '
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("E38").Value = '  -1.36%  '
$ws.Range("E39").Value = '  -1.47%  '
$ws.Range("E40").Value = '  -5.06%  '
$ws.Range("E41").Value = '  -0.35%  '
$ws.Range("E42").Value = '  -0.72%  '
$ws.Range("E43").Value = '  -1.57%  '
$ws.Range("E44").Value = '  +1.89%  '
$ws.Range("E45").Value = '  -0.45%  '
$ws.Range("E46").Value = '  -0.63%  '
$ws.Range("E47").Value = '  -2.49%  '
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("E49").Value = '  -1.09%  '
$ws.Range("E50").Value = '  -3.64%  '
$ws.Range("E51").Value = '  -2.41%  '
